# Update the Fruta/Hortaliza weekly price data.
# For each data row (2-16) the values in columns D, J, K, L, M, O, P are
# updated to reflect the latest weekly snapshot (effectively each row now
# carries the data that previously belonged to a different row in the
# historical series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Date (as Excel serial number), Volumen(J), PrecioMinimo(K), PrecioMaximo(L), PrecioPromedioPonderado(M), Origen(O), PrecioKg(P)
$data = @{
    2  = @{ D = 44987; J = 130; K = 4500; L = 5000; M = 4692; O = "Región Metropolitana";  P = 782 }
    3  = @{ D = 44957; J = 70;  K = 1500; L = 2000; M = 1857; O = "Región Metropolitana";  P = 310 }
    4  = @{ D = 45021; J = 50;  K = 4500; L = 5000; M = 4700; O = "Región Metropolitana";  P = 783 }
    5  = @{ D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";  P = 551 }
    6  = @{ D = 44876; J = 80;  K = 6500; L = 7000; M = 6812; O = "Región Metropolitana";  P = 1135 }
    7  = @{ D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";  P = 454 }
    8  = @{ D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    9  = @{ D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";  P = 484 }
    10 = @{ D = 44630; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";  P = 454 }
    11 = @{ D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";  P = 544 }
    12 = @{ D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
    13 = @{ D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";  P = 548 }
    14 = @{ D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";  P = 463 }
    15 = @{ D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";  P = 622 }
    16 = @{ D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";  P = 485 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O  # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals.P  # P - Precio $/Kg
}
